# "unify the conception of DataNode, DataTable, Entity."
#
# The sheet previously called "Property1" becomes "DataNode" (aligning its
# name with the DataNode/DataTable/Entity concepts used elsewhere), a couple
# of the data columns get a small manual width tweak, and the in-sheet
# selection is left on B41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: Property1 -> DataNode
$ws.Name = "DataNode"

# Small manual width adjustments on columns A and C
$ws.Columns.Item(1).ColumnWidth = 31.14
$ws.Columns.Item(3).ColumnWidth = 30.9

# Leave the cursor/selection on B41
[void]$ws.Range("B41").Select()
